$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (D,M,N,O,P,S)
$ws.Range("D2").Value = 44252
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 13000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 13500
$ws.Range("S2").Value = 750

# Row 3 (D,M)
$ws.Range("D3").Value = 44250
$ws.Range("M3").Value = 200

# Row 4 (D,M)
$ws.Range("D4").Value = 44253
$ws.Range("M4").Value = 160

# Row 5 (D,M,N,O,P,S)
$ws.Range("D5").Value = 44257
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("S5").Value = 806
